$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45 is new: fill constant columns copied from the product/market block,
# matching every other row in this sheet (same market, region, product, etc.).
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C45").Value = 'Arica y Parinacota'
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 'Fruta'
$ws.Range("G45").Value = 100108
$ws.Range("H45").Value = 'Tropicales y subtropicales'
$ws.Range("I45").Value = 100108001
$ws.Range("J45").Value = 'Guayaba'
$ws.Range("K45").Value = 'Sin especificar'
$ws.Range("Q45").Value = '$/kilo (en caja de 10 kilos )'
$ws.Range("R45").Value = 'Región de Arica y Parinacota'
$ws.Range("T45").Value = 1

# Weekly data refresh: each weeks row shifts down one slot (newest week
# inserted at row 14), carrying Fecha/Calidad/Volumen/Precio* along with it.

# Row 14
$ws.Range("D14").Value = 44764
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 500
$ws.Range("O14").Value = 600
$ws.Range("P14").Value = 550
$ws.Range("S14").Value = 550

# Row 15
$ws.Range("D15").Value = 44414
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 160
$ws.Range("N15").Value = 1300
$ws.Range("O15").Value = 1400
$ws.Range("P15").Value = 1350
$ws.Range("S15").Value = 1350

# Row 16
$ws.Range("D16").Value = 44750
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 700
$ws.Range("O16").Value = 800
$ws.Range("P16").Value = 750
$ws.Range("S16").Value = 750

# Row 17
$ws.Range("D17").Value = 44715
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 800
$ws.Range("O17").Value = 900
$ws.Range("P17").Value = 850
$ws.Range("S17").Value = 850

# Row 18
$ws.Range("D18").Value = 44715
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 600
$ws.Range("O18").Value = 700
$ws.Range("P18").Value = 650
$ws.Range("S18").Value = 650

# Row 19
$ws.Range("D19").Value = 44309
$ws.Range("D19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 160
$ws.Range("N19").Value = 1400
$ws.Range("O19").Value = 1500
$ws.Range("P19").Value = 1450
$ws.Range("S19").Value = 1450

# Row 20
$ws.Range("D20").Value = 44722
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 140
$ws.Range("N20").Value = 800
$ws.Range("O20").Value = 900
$ws.Range("P20").Value = 850
$ws.Range("S20").Value = 850

# Row 21
$ws.Range("D21").Value = 44722
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 700
$ws.Range("O21").Value = 800
$ws.Range("P21").Value = 750
$ws.Range("S21").Value = 750

# Row 22
$ws.Range("D22").Value = 44386
$ws.Range("D22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 160
$ws.Range("N22").Value = 700
$ws.Range("O22").Value = 750
$ws.Range("P22").Value = 725
$ws.Range("S22").Value = 725

# Row 23
$ws.Range("D23").Value = 44386
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 600
$ws.Range("O23").Value = 650
$ws.Range("P23").Value = 625
$ws.Range("S23").Value = 625

# Row 24
$ws.Range("D24").Value = 44407
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 600
$ws.Range("O24").Value = 650
$ws.Range("P24").Value = 625
$ws.Range("S24").Value = 625

# Row 25
$ws.Range("D25").Value = 44350
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 140
$ws.Range("N25").Value = 750
$ws.Range("O25").Value = 800
$ws.Range("P25").Value = 775
$ws.Range("S25").Value = 775

# Row 26
$ws.Range("D26").Value = 44498
$ws.Range("D26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L26").Value = 'Segunda'
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 1200
$ws.Range("O26").Value = 1300
$ws.Range("P26").Value = 1250
$ws.Range("S26").Value = 1250

# Row 27
$ws.Range("D27").Value = 44351
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 700
$ws.Range("O27").Value = 800
$ws.Range("P27").Value = 750
$ws.Range("S27").Value = 750

# Row 28
$ws.Range("D28").Value = 44351
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L28").Value = 'Segunda'
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 600
$ws.Range("O28").Value = 700
$ws.Range("P28").Value = 650
$ws.Range("S28").Value = 650

# Row 29
$ws.Range("D29").Value = 44687
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 120
$ws.Range("N29").Value = 1300
$ws.Range("O29").Value = 1400
$ws.Range("P29").Value = 1350
$ws.Range("S29").Value = 1350

# Row 30
$ws.Range("D30").Value = 44344
$ws.Range("D30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 140
$ws.Range("N30").Value = 1000
$ws.Range("O30").Value = 1200
$ws.Range("P30").Value = 1100
$ws.Range("S30").Value = 1100

# Row 31
$ws.Range("D31").Value = 44344
$ws.Range("D31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L31").Value = 'Segunda'
$ws.Range("M31").Value = 120
$ws.Range("N31").Value = 800
$ws.Range("O31").Value = 850
$ws.Range("P31").Value = 825
$ws.Range("S31").Value = 825

# Row 32
$ws.Range("D32").Value = 44260
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 1900
$ws.Range("O32").Value = 2000
$ws.Range("P32").Value = 1950
$ws.Range("S32").Value = 1950

# Row 33
$ws.Range("D33").Value = 44725
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 140
$ws.Range("N33").Value = 700
$ws.Range("O33").Value = 800
$ws.Range("P33").Value = 750
$ws.Range("S33").Value = 750

# Row 34
$ws.Range("D34").Value = 44725
$ws.Range("D34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L34").Value = 'Segunda'
$ws.Range("M34").Value = 160
$ws.Range("N34").Value = 500
$ws.Range("O34").Value = 600
$ws.Range("P34").Value = 550
$ws.Range("S34").Value = 550

# Row 35
$ws.Range("D35").Value = 44403
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 1200
$ws.Range("O35").Value = 1300
$ws.Range("P35").Value = 1250
$ws.Range("S35").Value = 1250

# Row 36
$ws.Range("D36").Value = 44403
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L36").Value = 'Segunda'
$ws.Range("M36").Value = 120
$ws.Range("N36").Value = 950
$ws.Range("O36").Value = 1000
$ws.Range("P36").Value = 975
$ws.Range("S36").Value = 975

# Row 37
$ws.Range("D37").Value = 44379
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 150
$ws.Range("N37").Value = 700
$ws.Range("O37").Value = 800
$ws.Range("P37").Value = 747
$ws.Range("S37").Value = 747

# Row 38
$ws.Range("D38").Value = 44379
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L38").Value = 'Segunda'
$ws.Range("M38").Value = 140
$ws.Range("N38").Value = 500
$ws.Range("O38").Value = 600
$ws.Range("P38").Value = 543
$ws.Range("S38").Value = 543

# Row 39
$ws.Range("D39").Value = 44348
$ws.Range("D39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 120
$ws.Range("N39").Value = 1000
$ws.Range("O39").Value = 1100
$ws.Range("P39").Value = 1050
$ws.Range("S39").Value = 1050

# Row 40
$ws.Range("D40").Value = 44389
$ws.Range("D40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 140
$ws.Range("N40").Value = 750
$ws.Range("O40").Value = 800
$ws.Range("P40").Value = 775
$ws.Range("S40").Value = 775

# Row 41
$ws.Range("D41").Value = 44389
$ws.Range("D41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L41").Value = 'Segunda'
$ws.Range("M41").Value = 120
$ws.Range("N41").Value = 600
$ws.Range("O41").Value = 700
$ws.Range("P41").Value = 650
$ws.Range("S41").Value = 650

# Row 42
$ws.Range("D42").Value = 44417
$ws.Range("D42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L42").Value = 'Primera'
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = 1300
$ws.Range("O42").Value = 1400
$ws.Range("P42").Value = 1350
$ws.Range("S42").Value = 1350

# Row 43
$ws.Range("D43").Value = 44425
$ws.Range("D43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L43").Value = 'Primera'
$ws.Range("M43").Value = 140
$ws.Range("N43").Value = 1200
$ws.Range("O43").Value = 1300
$ws.Range("P43").Value = 1250
$ws.Range("S43").Value = 1250

# Row 44
$ws.Range("D44").Value = 44358
$ws.Range("D44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L44").Value = 'Primera'
$ws.Range("M44").Value = 200
$ws.Range("N44").Value = 700
$ws.Range("O44").Value = 800
$ws.Range("P44").Value = 750
$ws.Range("S44").Value = 750

# Row 45
$ws.Range("D45").Value = 44358
$ws.Range("D45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L45").Value = 'Segunda'
$ws.Range("M45").Value = 200
$ws.Range("N45").Value = 600
$ws.Range("O45").Value = 650
$ws.Range("P45").Value = 625
$ws.Range("S45").Value = 625
